# Expand dropdown providers with comprehensive value ranges.
#
# The "Default" (column D) cells for each Slider option on the "Menu Mock"
# sheet were storing their default as formatted display text (e.g. "25%",
# "0.28x", "2.5s") that duplicated the adjacent "Values" (column E) dropdown
# list. Converting D to a plain numeric value (e.g. 0.25, 0.28, 2.5) lets the
# preset sync pick a value that actually matches an entry in the expanded
# provider list, instead of every trigger collapsing onto the same value.
#
# The old column E text (the literal list of allowed values, duplicated
# verbatim on the "Providers" sheet in column B) is cleared here because the
# provider is being expanded to a much larger generated range elsewhere;
# the cooldown rows keep a single "0s" placeholder since 0 remains a valid
# cooldown value across all triggers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

# CategoryCustomBasic
$ws.Range("D42").Value = 0.25   # Chance:    25%   -> 0.25
$ws.Range("E42").Value = ""
$ws.Range("D43").Value = 0.28   # TimeScale: 0.28x -> 0.28
$ws.Range("E43").Value = ""
$ws.Range("D44").Value = 2.5    # Duration:  2.5s  -> 2.5
$ws.Range("E44").Value = ""
$ws.Range("D45").Value = 10     # Cooldown:  10.0s -> 10
$ws.Range("E45").Value = "0s"

# CategoryCustomCritical
$ws.Range("D51").Value = 0.75   # Chance:    75%   -> 0.75
$ws.Range("E51").Value = ""
$ws.Range("D52").Value = 0.25   # TimeScale: 0.25x -> 0.25
$ws.Range("E52").Value = ""
$ws.Range("D53").Value = 3.0    # Duration:  3.0s  -> 3.0
$ws.Range("E53").Value = ""
$ws.Range("D54").Value = 10     # Cooldown:  10.0s -> 10
$ws.Range("E54").Value = "0s"

# CategoryCustomDismemberment
$ws.Range("D60").Value = 0.6    # Chance:    60%   -> 0.6
$ws.Range("E60").Value = ""
$ws.Range("D61").Value = 0.3    # TimeScale: 0.30x -> 0.3
$ws.Range("E61").Value = ""
$ws.Range("D62").Value = 2.0    # Duration:  2.0s  -> 2.0
$ws.Range("E62").Value = ""
$ws.Range("D63").Value = 10     # Cooldown:  10.0s -> 10
$ws.Range("E63").Value = "0s"

# CategoryCustomDismemberment-adjacent block (90% chance trigger)
$ws.Range("D69").Value = 0.9    # Chance:    90%   -> 0.9
$ws.Range("E69").Value = ""
$ws.Range("D70").Value = 0.23   # TimeScale: 0.23x -> 0.23
$ws.Range("E70").Value = ""
$ws.Range("D71").Value = 3.25   # Duration:  3.25s -> 3.25
$ws.Range("E71").Value = ""
$ws.Range("D72").Value = 10     # Cooldown:  10.0s -> 10
$ws.Range("E72").Value = "0s"

# CategoryCustomLastEnemy
$ws.Range("D78").Value = 1.0    # Chance:    100%  -> 1.0
$ws.Range("E78").Value = ""
$ws.Range("D79").Value = 0.26   # TimeScale: 0.26x -> 0.26
$ws.Range("E79").Value = ""
$ws.Range("D80").Value = 2.75   # Duration:  2.75s -> 2.75
$ws.Range("E80").Value = ""
$ws.Range("D81").Value = 20     # Cooldown:  20.0s -> 20
$ws.Range("E81").Value = "0s"

# CategoryCustomLastStand
$ws.Range("D87").Value = 0.21   # TimeScale: 0.21x -> 0.21
$ws.Range("E87").Value = ""
$ws.Range("D88").Value = 3.5    # Duration:  3.5s  -> 3.5
$ws.Range("E88").Value = ""
$ws.Range("D89").Value = 60     # Cooldown:  60.0s -> 60
$ws.Range("E89").Value = "0s"

# CategoryCustomParry
$ws.Range("D94").Value = 0.5    # Chance:    50%   -> 0.5
$ws.Range("E94").Value = ""
$ws.Range("D95").Value = 0.34   # TimeScale: 0.34x -> 0.34
$ws.Range("E95").Value = ""
$ws.Range("D96").Value = 1.5    # Duration:  1.5s  -> 1.5
$ws.Range("E96").Value = ""
$ws.Range("D97").Value = 5      # Cooldown:  5.0s  -> 5
$ws.Range("E97").Value = "0s"

# "Providers" reference sheet: same value lists duplicated in column B next
# to each provider name in column A. These collapse the same way as the
# column E lists above (CustomCooldownProvider keeps its "0s" placeholder).
$wsProv = $wb.Worksheets.Item("Providers")
$wsProv.Range("B9").Value = ""     # CustomChanceProvider
$wsProv.Range("B10").Value = "0s"  # CustomCooldownProvider
$wsProv.Range("B11").Value = ""    # CustomDurationProvider
$wsProv.Range("B14").Value = ""    # CustomTimeScaleProvider

Write-Host "Expanded dropdown providers: cleared stale Values lists and normalized Default cells to numeric."
